# Update "想去人数" (F column) counts on sheets 展览, 演出, 全部类型
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 286
$ws1.Range("F6").Value = 99
$ws1.Range("F8").Value = 217
$ws1.Range("F9").Value = 2020
$ws1.Range("F11").Value = 4791
$ws1.Range("F12").Value = 88

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 8
$ws2.Range("F5").Value = 13

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 286
$ws4.Range("F8").Value = 99
$ws4.Range("F10").Value = 217
$ws4.Range("F11").Value = 8
$ws4.Range("F12").Value = 13
$ws4.Range("F13").Value = 2020
$ws4.Range("F15").Value = 4791
$ws4.Range("F16").Value = 88
